$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-04 Saturday" "2025-01-05 Sunday"

Replace-Text "802÷2=" "812÷4="
Replace-Text "127÷8=" "681÷9="
Replace-Text "964÷3=" "123÷5="
Replace-Text "748÷4=" "896÷6="
Replace-Text "355÷9=" "156÷6="

Replace-Text "281÷7=" "716÷5="
Replace-Text "370÷8=" "217÷9="
Replace-Text "851÷5=" "285÷3="
Replace-Text "640÷7=" "841÷2="
Replace-Text "435÷8=" "437÷4="

Replace-Text "210÷9=" "418÷5="
Replace-Text "298÷5=" "533÷9="
Replace-Text "247÷6=" "832÷9="
Replace-Text "535÷7=" "832÷4="
Replace-Text "129÷7=" "341÷8="

Replace-Text "217÷2=" "198÷3="
Replace-Text "942÷7=" "410÷2="
Replace-Text "578÷8=" "288÷9="
Replace-Text "555÷9=" "281÷6="
Replace-Text "828÷6=" "484÷8="

Replace-Text "711÷9=" "847÷4="
Replace-Text "244÷9=" "941÷6="
Replace-Text "303÷6=" "725÷5="
Replace-Text "768÷6=" "740÷9="
Replace-Text "969÷2=" "736÷6="

Write-Output "Done"
